$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew
$overview.Range("E3").Value = $statusNew
$overview.Range("F3").Value = $statusNew

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusNew
$zhcn.Range("C3").Value = $statusNew

$zhcn.Range("I2").Value = "1b1858c3-3b00-47a2-8f24-973e2b5359c5.md"
$zhcn.Range("J2").Value = "1b1858c3-3b00-47a2-8f24-973e2b5359c5.12bba9edd97355a63be4dcec563f4f1a1f89b5b0.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-02 11:12:50"

$zhcn.Range("I3").Value = "c418e212-b76f-4bce-8821-99186b269b0e.md"
$zhcn.Range("J3").Value = "c418e212-b76f-4bce-8821-99186b269b0e.5917ff04743a22e1bdeea04a95e466e37b6c7d73.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-02 11:12:50"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90b0e6d9aae3728d2d5a708251d9b1755ac564a5/e2e/1b1858c3-3b00-47a2-8f24-973e2b5359c5.md", "", "", "1b1858c3-3b00-47a2-8f24-973e2b5359c5.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90b0e6d9aae3728d2d5a708251d9b1755ac564a5/e2e/c418e212-b76f-4bce-8821-99186b269b0e.md", "", "", "c418e212-b76f-4bce-8821-99186b269b0e.md")

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusNew
$dede.Range("C3").Value = $statusNew

$dede.Range("I2").Value = "1b1858c3-3b00-47a2-8f24-973e2b5359c5.md"
$dede.Range("J2").Value = "1b1858c3-3b00-47a2-8f24-973e2b5359c5.12bba9edd97355a63be4dcec563f4f1a1f89b5b0.de-de.xlf"
$dede.Range("K2").Value = "2016-09-02 11:12:58"

$dede.Range("I3").Value = "c418e212-b76f-4bce-8821-99186b269b0e.md"
$dede.Range("J3").Value = "c418e212-b76f-4bce-8821-99186b269b0e.5917ff04743a22e1bdeea04a95e466e37b6c7d73.de-de.xlf"
$dede.Range("K3").Value = "2016-09-02 11:12:58"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90b0e6d9aae3728d2d5a708251d9b1755ac564a5/e2e/1b1858c3-3b00-47a2-8f24-973e2b5359c5.md", "", "", "1b1858c3-3b00-47a2-8f24-973e2b5359c5.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90b0e6d9aae3728d2d5a708251d9b1755ac564a5/e2e/c418e212-b76f-4bce-8821-99186b269b0e.md", "", "", "c418e212-b76f-4bce-8821-99186b269b0e.md")
